$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G5").Value = 'AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G7").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G9").Value = 'yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg'
$ws.Range("G12").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G18").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G19").Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G20").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G21").Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G24").Value = 'AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G26").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G28").Value = 'dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range("G31").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G37").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G38").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G39").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G40").Value = 'alshimaa.atef@med.asu.edu.egm, shaimaa.ahmed@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G43").Value = 'Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G45").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G47").Value = 'esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range("G50").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G56").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G57").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G58").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G59").Value = 'nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G63").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G66").Value = 'dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range("G75").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G76").Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G77").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G78").Value = 'nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G81").Value = 'AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg'
$ws.Range("G85").Value = 'dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range("G88").Value = 'heba.e.mohamed@med.asu.edu.eg, mona.I.hussein@med.asu.edu.eg, dalia.t.abdelmegid@med.asu.edu.eg'
$ws.Range("G94").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G95").Value = 'marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G96").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G97").Value = 'nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, heba@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range("G100").Value = 'Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G101").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G102").Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G104").Value = 'esraa.mostafa@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg'
$ws.Range("G113").Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G115").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
